# Update Seb's Booked Hours timesheet:
#  - Fill in row 11 (Wk2, 29/07/21) "Work on improving the github setup..." session
#  - Fill in row 12 (Wk2, 01/08/21) "Sent email to auditor asking to organise a meeting" session
#  - These were accidentally left blank/under-recorded; entering the real meeting time.
#  - Incidental UI state: the author ended up with cell L3 selected on each sheet
#    (an accidental extra click that created an empty styled cell on BookedHours).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

# --- Row 11: 29/07/2021, 18:10 -> 22:24 ---
$ws1.Range("A11").Value = 2
$ws1.Range("B11").Value = 44405
$ws1.Range("C11").Value = 0.756944444444444
$ws1.Range("D11").Value = 44405
$ws1.Range("E11").Value = 0.933333333333333
$ws1.Range("F11").Value = "Work on improving the github setup – added project/task board and issues, and switched to an organsiation instead of personal as that gives us better project management features"

# --- Row 12: 01/08/2021, 10:20 -> 10:25 ---
$ws1.Range("A12").Value = 2
$ws1.Range("B12").Value = 44408
$ws1.Range("C12").Value = 0.430555555555556
$ws1.Range("D12").Value = 44408
$ws1.Range("E12").Value = 0.434027777777778
$ws1.Range("F12").Value = "Sent email to auditor asking to organise a meeting"

# --- Incidental: an extra cell (L3) got touched/selected on BookedHours ---
$ws1.Range("L3").NumberFormat = "General"
$ws1.Range("L3").Select()

# --- Incidental: selection state left behind on the other two sheets ---
$ws2.Activate()
$ws2.Range("D22").Select()

$ws3.Activate()
$ws3.Range("A1").Select()

# Restore BookedHours as the active/selected sheet (it was the sheet in view)
$ws1.Activate()
